$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$rows = @(2,3,4,5,6,7)
$vals = @(4714,5082,1254,299,4316,15665)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Logan Square')
$rows = @(2,6,7)
$vals = @(53,61,179)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Austin')
$rows = @(2,5,7)
$vals = @(300,38,1040)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Garfield Park')
$rows = @(3,7)
$vals = @(242,713)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('New City')
$rows = @(2,7)
$vals = @(108,304)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Fuller Park')
$rows = @(6,7)
$vals = @(32,74)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('By Neighborhood')
$rows = @(6,7,8,10,18,19,20,21,22,27,29,30,33,42,53,54,55,58,63,64,65,67,70,72,77,85,88,89,91,94,96,97,101)
$vals = @(119,513,1040,104,110,430,394,48,46,139,863,74,713,512,179,329,150,9,41,109,304,539,27,61,105,804,172,225,211,194,174,135,15665)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('North Lawndale')
$rows = @(3,7)
$vals = @(210,539)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Loop')
$rows = @(2,4,7)
$vals = @(60,28,329)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Englewood')
$rows = @(3,6,7)
$vals = @(326,225,863)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Chatham')
$rows = @(3,5,6,7)
$vals = @(134,7,123,430)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Ashburn')
$rows = @(6,7)
$vals = @(22,119)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Humboldt Park')
$rows = @(6,7)
$vals = @(144,512)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Avondale')
$rows = @(3,7)
$vals = @(28,104)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Lower West Side')
$rows = @(4,7)
$vals = @(10,150)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('West Ridge')
$rows = @(3,7)
$vals = @(48,174)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Washington Park')
$rows = @(3,7)
$vals = @(94,211)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Chinatown')
$rows = @(6,7)
$vals = @(25,48)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Near South Side')
$rows = @(3,6,7)
$vals = @(28,31,109)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Chicago Lawn')
$rows = @(2,7)
$vals = @(120,394)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Calumet Heights')
$rows = @(2,7)
$vals = @(40,110)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Auburn Gresham')
$rows = @(6,7)
$vals = @(125,513)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('West Loop')
$rows = @(3,7)
$vals = @(44,194)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('West Town')
$rows = @(3,6,7)
$vals = @(30,67,135)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('O''Hare')
$rows = @(5,7)
$vals = @(2,27)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('United Center')
$rows = @(2,3,7)
$vals = @(48,64,172)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Uptown')
$rows = @(2,7)
$vals = @(63,225)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Edgewater')
$rows = @(2,6,7)
$vals = @(38,40,139)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('South Shore')
$rows = @(2,4,7)
$vals = @(241,50,804)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Clearing')
$rows = @(3,7)
$vals = @(17,46)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Old Town')
$rows = @(3,7)
$vals = @(15,61)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Riverdale')
$rows = @(2,7)
$vals = @(36,105)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}

$ws = $wb.Worksheets.Item('Millenium Park')
$rows = @(6,7)
$vals = @(3,9)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 12).Value = $vals[$i]
}
